# Add a new Eurobarometer wave entry (ZA7782 / EB 95.2) as a new row
# inserted above the current row 3, pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 3 (inherits formatting, e.g. column B's
# text/quote-prefix style, from the row below it).
$ws.Rows.Item(3).EntireRow.Insert() | Out-Null

# Fill in the new survey metadata row.
$ws.Range("A3").Value = "ZA7782"
# Leading apostrophe keeps "95.2" stored as text (matches the existing
# quotePrefix style already applied to this cell by the row insert).
$ws.Range("B3").Value = "'95.2"
$ws.Range("C3").Value = "April-May 2021"
$ws.Range("D3").Value = "European citizens" + [char]0x2019 + " knowledge and attitudes towards science and technology"

# Match the author's final cursor position after adding the row.
$ws.Range("D4").Select() | Out-Null
